$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 78 (shifts old rows 78-82 down to 79-83)
$ws.Rows.Item(78).Insert()

# Copy style of date cell (D79, which was old D78) into new D78
$ws.Range("D79").Copy()
$ws.Range("D78").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 45075
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112001
$ws.Range("G78").Value = "Berenjena"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 40
$ws.Range("K78").Value = 7000
$ws.Range("L78").Value = 8000
$ws.Range("M78").Value = 7500
$ws.Range("N78").Value = "$/caja 60 unidades"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 125
$ws.Range("Q78").Value = 60
$ws.Range("R78").Value = "Hortaliza"
